$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "< 10 Manuf., Agriculture, and Forestry`n< 5 Wholesale, and Retail`n< 10 Services, and Mining"
$ws.Range("C10").Value = "< CFA 20 Millionlion Manuf., Agriculture, and Forestry`n< CFA 15 Millionlion Wholesale`n< CFA 10 Millionlion Retail, Services, and Mining"

$ws.Range("B11").Value = "< 50 Manuf., Agriculture, Forestry, and Wholesale`n< 30 Retail, Services, and Mining"
$ws.Range("C11").Value = "< CFA 250 Millionlion Manuf., Agriculture, and Forestry `n< CFA 150 Millionlion Wholesale`n< CFA 50 Millionlion Retail`n< CFA 75 Millionlion Services, and Mining"
$ws.Range("D11").Value = "< CFA 250 Millionlion Manuf., Agriculture, and Forestry `n< CFA 200 Millionlion Wholesale`n< CFA 100 Millionlion Retail, Services, and Mining"

$ws.Range("B12").Value = "< 100 Manuf., Agriculture, and Forestry`n< 50 Wholesale, Retail, Services, and Mining"
$ws.Range("C12").Value = "< CFA 750 Millionlion Manuf., Agriculture, and Forestry `n< CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"
$ws.Range("D12").Value = "< CFA 500 Millionlion Manuf., Agriculture, and Forestry `n< CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"

$ws.Range("B13").Value = "> 100 Manuf., Agriculture, and Forestry`n> 50 Wholesale, Retail, Services, and Mining"
$ws.Range("C13").Value = "> CFA 750 Millionlion Manuf., Agriculture, and Forestry `n> CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"
$ws.Range("D13").Value = "> CFA 500 Millionlion Manuf., Agriculture, and Forestry `n> CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"
